$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per latest scrape
$ws.Range("D2").Value = "68.104.50"
$ws.Range("D3").Value = "3.587.56"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +10.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "568.56"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").Value = "3.583.61"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.43"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +8.36%  "
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000280"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.23"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.45%  "
$ws.Range("D15").Value = "4.158.26"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "3.590.15"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.89"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.54%  "
$ws.Range("D19").Value = "67.943.81"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "401.41"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.12%  "
$ws.Range("B23").Value = "RenderToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.52"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +11.93%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.18"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.40"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.88"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.43"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.88"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +9.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.21"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.74"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.46"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "669.21"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.06"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.113"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.35"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "41.14"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.407"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.26"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +11.35%  "
$ws.Range("D40").Value = "0.0₃0752"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").Value = "3.191.39"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.80"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +12.32%  "
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.130"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.69"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.05"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.60"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +10.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "139.03"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.15%  "
